$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accelerometer readings (x, y, z), one new sample added near the top
# of the data and the data shifted down by one row; the sheet's last row
# (old row 22) drops off the bottom so the sheet shrinks from A1:C22 to
# A1:C21.
$data = @(
    @(-3.540287351608276,  5.382533311843872,  -1.780441856384277),
    @(-3.185474908351898,  5.544018769264222,  -1.341128690540791),
    @(-3.169438767433166,  5.001180601119994,  -0.6735752552747722),
    @(-3.514218628406524,  5.064400196075439,   0.5307594910264006),
    @(-3.949312973022461,  3.677936935424804,   0.8671627283096315),
    @(-4.74082317352295,   2.939013087749481,   1.662665629386903),
    @(-6.156321334838866,  2.482692444324494,   2.901540523767471),
    @(-6.70512387752533,   1.984266856312751,   2.305714881420135),
    @(-6.514288425445558,  2.538958132266995,   0.7396240234375031),
    @(13.64799528121941,  -0.5907153248786794, -5.010972028970702),
    @(25.54574513435358,  -2.183045053482042,  -7.100710201263412),
    @(-1.683456826210012,  4.46893746256828,    1.391087603569029),
    @(1.132775115966798,   4.239795589447025,   0.7271292686462392),
    @(1.974032163619995,   6.137257993221282,   0.4547893404960567),
    @(1.989514970779418,   6.169636392593385,   2.127831518650056),
    @(2.919293570518496,   4.336295771598811,  -1.078460484743126),
    @(3.384284400939942,   5.661798715591432,   3.016176247596746),
    @(1.664549851417542,   6.045073473453522,   1.692821365594864),
    @(1.458636522293091,   6.004646182060242,   0.9704791456460963),
    @(1.448996758460998,   5.705842137336731,   1.135319881141186)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
    $ws.Range("C$row").Value = $data[$i][2]
}

# Remove the now-stale final row (old row 22) so the used range shrinks
# back down to A1:C21.
$ws.Range("A22:C22").ClearContents()
